$wb = $excel.ActiveWorkbook

# --- "Monthly Costs" sheet: add two new data rows (2 and 3) ---
$wsCosts = $wb.Worksheets.Item("Monthly Costs")

# Text-typed columns must keep their literal text (e.g. dates as plain
# strings, not auto-converted to date serials), so force a text format
# before writing the values.
$wsCosts.Range("A2:C3").NumberFormat = "@"

$wsCosts.Range("A2").Value = "2024-10-10"
$wsCosts.Range("B2").Value = "Example Company"
$wsCosts.Range("C2").Value = "Item 1"
$wsCosts.Range("D2").Value = 5
$wsCosts.Range("E2").Value = 100
$wsCosts.Range("F2").Value = 500

$wsCosts.Range("A3").Value = "2024-10-10"
$wsCosts.Range("B3").Value = "Example Company"
$wsCosts.Range("C3").Value = "Item 2"
$wsCosts.Range("D3").Value = 3
$wsCosts.Range("E3").Value = 200
$wsCosts.Range("F3").Value = 600

# --- "Pricing" sheet: add two new data rows (2 and 3) ---
$wsPricing = $wb.Worksheets.Item("Pricing")

$wsPricing.Range("A2:A3").NumberFormat = "@"

$wsPricing.Range("A2").Value = "Item 1"
$wsPricing.Range("B2").Value = 100

$wsPricing.Range("A3").Value = "Item 2"
$wsPricing.Range("B3").Value = 200
